$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: becomes the old Row 16 record (Garnlav @ 570817/6736417, reporter Bo karlstens) ---
$ws.Range("A15").Value2 = 130983071
$ws.Range("B15").Value2 = 79244
$ws.Range("E15").Value2 = 6425
$ws.Range("F15").Value2 = "Garnlav"
$ws.Range("G15").Value2 = "Alectoria sarmentosa"
$ws.Range("H15").Value2 = "(Ach.) Ach."
$ws.Range("Q15").Value2 = 570817
$ws.Range("R15").Value2 = 6736417
$ws.Range("S15").Value2 = 10
$ws.Range("Z15").Value2 = "08:53"
$ws.Range("AB15").Value2 = "08:53"
$ws.Range("AC15").ClearContents()
$ws.Range("AF15").Value2 = ""
$ws.Range("AW15").Value2 = "Bo karlstens"
$ws.Range("AX15").Value2 = "Bo karlstens"

# --- Row 16: becomes the old Row 18 record (Garnlav @ 570825/6736389, reporter Göran Ehn) ---
$ws.Range("A16").Value2 = 130983619
$ws.Range("P16").Value2 = "Flytjärnsmyren, Dlr"
$ws.Range("Q16").Value2 = 570825
$ws.Range("R16").Value2 = 6736389
$ws.Range("Z16").Value2 = "08:54"
$ws.Range("AB16").Value2 = "08:54"
$ws.Range("AF16").ClearContents()
$ws.Range("AW16").Value2 = "Göran Ehn"
$ws.Range("AX16").Value2 = "Göran Ehn"

# --- Row 18: becomes the old Row 15 record (Tretåig hackspett @ 570952/6736563, reporter Erik Danielsson) ---
$ws.Range("A18").Value2 = 130979082
$ws.Range("B18").Value2 = 57884
$ws.Range("E18").Value2 = 100109
$ws.Range("F18").Value2 = "Tretåig hackspett"
$ws.Range("G18").Value2 = "Picoides tridactylus"
$ws.Range("H18").Value2 = "(Linnaeus, 1758)"
$ws.Range("P18").Value2 = "Flytjärnsmyran, Dlr"
$ws.Range("Q18").Value2 = 570952
$ws.Range("R18").Value2 = 6736563
$ws.Range("S18").Value2 = 1
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
$ws.Range("AC18").Value2 = "Äldre ringhack"
$ws.Range("AW18").Value2 = "Erik Danielsson"
$ws.Range("AX18").Value2 = "Erik Danielsson"
